$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FatosIn")
$ws.Cells.Item(28,2).Value = "XYZ_UNIQUE_TEST"
